$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''309.35'
$ws.Range('E2').Value = '''-0.84%'
$ws.Range('D3').Value = '''37.02'
$ws.Range('E3').Value = '''-2.13%'
$ws.Range('D4').Value = '''5.131'
$ws.Range('E4').Value = '''-0.04%'
$ws.Range('D5').Value = '''0.07847'
$ws.Range('E5').Value = '''-0.97%'
$ws.Range('B6').Value = '''GateToken'
$ws.Range('C6').Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').Value = '''4.403'
$ws.Range('E6').Value = '''0.02%'
$ws.Range('B7').Value = '''KuCoinToken'
$ws.Range('C7').Value = '''https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '''8.271'
$ws.Range('E7').Value = '''0.28%'
$ws.Range('B8').Value = '''FTXToken'
$ws.Range('C8').Value = '''https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').Value = '''1.877'
$ws.Range('E8').Value = '''-1.68%'
$ws.Range('B9').Value = '''BTSEToken'
$ws.Range('C9').Value = '''https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').Value = '''2.982'
$ws.Range('E9').Value = '''4.91%'
$ws.Range('B10').Value = '''MXToken'
$ws.Range('C10').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D10').Value = '''0.9240'
$ws.Range('E10').Value = '''-0.36%'
$ws.Range('B11').Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '''0.1171'
$ws.Range('E11').Value = '''-2.56%'
$ws.Range('B12').Value = '''WazirX'
$ws.Range('C12').Value = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D12').Value = '''0.1900'
$ws.Range('E12').Value = '''-0.48%'
$ws.Range('B13').Value = '''MandalaExchangeToken'
$ws.Range('C13').Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').Value = '''0.08911'
$ws.Range('E13').Value = '''-4.10%'
$ws.Range('B14').Value = '''BitrueCoin'
$ws.Range('C14').Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = '''0.03318'
$ws.Range('E14').Value = '''-1.88%'
$ws.Range('B15').Value = '''BitMartToken'
$ws.Range('C15').Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '''0.09606'
$ws.Range('E15').Value = '''-0.16%'
$ws.Range('B16').Value = '''BitForexToken'
$ws.Range('C16').Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '''0.001382'
$ws.Range('E16').Value = '''0.67%'
$ws.Range('B17').Value = '''TigerCash'
$ws.Range('C17').Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = '''0.006202'
$ws.Range('E17').Value = '''6.57%'
$ws.Range('B18').Value = '''LEO'
$ws.Range('C18').Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = '''3.394'
$ws.Range('E18').Value = '''-3.85%'
$ws.Range('D19').Value = '''0.3458'
$ws.Range('E19').Value = '''0.25%'
$ws.Range('D20').Value = '''6.367'
$ws.Range('E20').Value = '''20.98%'
$ws.Range('D21').Value = '''0.1297'
$ws.Range('E21').Value = '''1.12%'
$ws.Range('D22').Value = '''0.2406'
$ws.Range('E22').Value = '''-6.95%'
$ws.Range('D23').Value = '''0.04343'
$ws.Range('E23').Value = '''-0.43%'
$ws.Range('D24').Value = '''0.001202'
$ws.Range('E24').Value = '''-3.82%'
$ws.Range('D25').Value = '''0.004281'
$ws.Range('D26').Value = '''0.0001401'
$ws.Range('E26').Value = '''8.16%'
$ws.Range('D27').Value = '''0.0002902'
$ws.Range('D39').Value = '''0.02161'
$ws.Range('E39').Value = '''3.17%'
$ws.Range('D40').Value = '''0.05005'
$ws.Range('E40').Value = '''-1.34%'
$ws.Range('D41').Value = '''0.007590'
$ws.Range('E41').Value = '''-0.50%'
$ws.Range('E42').Value = '''0.19%'
$ws.Range('E43').Value = '''-6.60%'
$ws.Range('D44').Value = '''0.002009'
$ws.Range('E44').Value = '''-0.27%'
$ws.Range('D45').Value = '''0.008912'
$ws.Range('E45').Value = '''3.06%'
$ws.Range('D46').Value = '''0.00006578'
$ws.Range('E46').Value = '''-1.84%'
$ws.Range('D47').Value = '''0.00000000751'
$ws.Range('E47').Value = '''0.13%'
$ws.Range('D48').Value = '''0.003296'
$ws.Range('E48').Value = '''13.69%'
$ws.Range('D49').Value = '''0.001444'
$ws.Range('E49').Value = '''20.39%'
$ws.Range('D50').Value = '''0.00002102'
$ws.Range('E50').Value = '''0.13%'
$ws.Range('D51').Value = '''0.0002002'
$ws.Range('E51').Value = '''0.13%'
